# Apply crypto price/volume updates to Sheet1 (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.000" -> 1, "277.20" -> 277.2) must be forced to Text format
# first so the literal digit string is preserved, matching the source data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"

# --- Column D (Price) and Column E (Volume(1h)) updates ---
$ws.Range("D2").Value = "26.818.35"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "1.875.94"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "277.20"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5282"
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("D8").Value = "0.3415"
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").Value = "0.06943"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("D10").Value = "20.05"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "0.8037"
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("D12").Value = "0.07712"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.876.29"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "5.190"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").Value = "90.20"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "0.9991"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "0.000008050"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D20").Value = "26.845.34"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").Value = "2.086.43"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "4.751"
$ws.Range("D23").Value = "10.03"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "6.166"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "2.392"
$ws.Range("E25").Value = "  +8.64%  "
$ws.Range("D26").Value = "146.54"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("D27").Value = "17.34"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "1.656"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "113.69"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "4.307"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "0.08911"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "0.04925"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").Value = "0.7250"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "2.871"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "3.277"
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").Value = "2.343"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "0.01857"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "0.5134"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "0.9548"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "116.08"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("D43").Value = "6.152"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "8.111"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "0.4468"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").Value = "9.308"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "36.26"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "1.491"
$ws.Range("E51").Value = "  -0.55%  "
